$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '86.910.01'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +9.33%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '3.319.06'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +4.93%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.00'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.12%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '218.71'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +5.96%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '638.53'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +1.88%  '
$ws.Range('E7').Value = '  +20.44%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.998'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  -0.11%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.614'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +3.97%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '3.321.02'
$ws.Range('D10').Style = 'Normal'
$ws.Range('E10').Value = '  +5.09%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.604'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  +1.45%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '0.0000275'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  +9.53%  '
$ws.Range('E13').Value = '  +1.70%  '
$ws.Range('B14').Value = 'WrappedliquidstakedEther2.0'
$ws.Range('C14').Value = 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '3.931.85'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +5.17%  '
$ws.Range('B15').Value = 'Avalanche'
$ws.Range('C15').Value = 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '34.30'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +8.96%  '
$ws.Range('E16').Value = '  +2.24%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '86.801.84'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +9.05%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '3.311.78'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +4.82%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '14.69'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  +2.53%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '3.20'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +8.56%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '447.45'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +2.28%  '
$ws.Range('E22').Value = '  +0.48%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '5.25'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '7.45'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +9.92%  '
$ws.Range('E25').Value = '  +15.32%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '12.24'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +12.48%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '3.436.38'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +3.28%  '
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '78.57'
$ws.Range('D28').Style = 'Normal'
$ws.Range('E28').Value = '  +3.08%  '
$ws.Range('E29').Value = '  +8.67%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.999'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  -0.01%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '0.174'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  +45.08%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '606.96'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  +10.12%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '9.29'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +3.17%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '0.999'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +0.18%  '
$ws.Range('E35').Value = '  +5.60%  '
$ws.Range('E36').Value = '  +2.50%  '
$ws.Range('E37').Value = '  +2.16%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '23.32'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.83%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '6.50'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +15.48%  '
$ws.Range('E40').Value = '  +2.98%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '0.998'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  -0.06%  '
$ws.Range('B42').Value = 'dogwifhat'
$ws.Range('C42').Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '3.16'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +17.79%  '
$ws.Range('B43').Value = 'WhiteBITCoin'
$ws.Range('C43').Value = 'https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '21.32'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +2.61%  '
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '2.06'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +14.40%  '
$ws.Range('B45').Value = 'USDe'
$ws.Range('C45').Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '1.00'
$ws.Range('D45').Style = 'Normal'
$ws.Range('E45').Value = '  +0.07%  '
$ws.Range('B46').Value = 'Monero'
$ws.Range('C46').Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '156.50'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -4.16%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '189.40'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  +0.80%  '
$ws.Range('E48').Value = '  +5.35%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '1.36'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +4.85%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.784'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +0.48%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '26.49'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +6.52%  '
